$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# fix data text metric names: strip stray leading/trailing whitespace and
# space-separated casing from the "Self Employed" / "Unemployed" / "Inactive"
# metric labels in column A so they match the tidy metric-code convention
# used elsewhere in the column (e.g. empRate, unempRate, inactiveRate).
$ws.Range("A7").Value = "SelfEmployed"
$ws.Range("A8").Value = "Unemployed"
$ws.Range("A9").Value = "Inactive"

# restore the window/selection state that was saved with the workbook
# (scrolled so row 5 is at the top, with E6 the active/selected cell)
$ws.Activate()
$ws.Range("E6").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
